$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Rows 30 and 31: mark as complete (previously "in progress")
$ws.Range("H30").Value = "complete"
$ws.Range("I30").Value = "complete"
$ws.Range("H31").Value = "complete"
$ws.Range("I31").Value = "complete"

# Row 32: new reinsurance test entry (fac separated from surplus share, ri1a added for fac)
# Set E32 before C32 so the shared-string table gets the two new strings in the
# same order as the target file (Calcrule string first, Description string second).
$ws.Range("E32").Value = "14,21, 23,25"
$ws.Range("C32").Value = "Reinsurance example with location level fac and surplus share on a subset of locations and two per risk treaties on all locations"
$ws.Range("D32").Value = "0,2"

# Copy number/alignment formatting from the row above for F32:G32 (right-aligned style)
$ws.Range("F31:G31").Copy()
$ws.Range("F32:G32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 2

$ws.Range("H32").Value = "in progress"
$ws.Range("I32").Value = "in progress"
